$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move SUB/SHR from row 28 into row 27 (C27/E27), and clear row 28's old content
$ws.Range("C27").Value = "SUB"
$ws.Range("E27").Value = "SHR"

# Row 28 becomes just "PAR" in B28; clear old D28 entirely
$ws.Range("D28").ClearContents()
$ws.Range("B28").Value = "PAR"

# Update the active selection to E27
$ws.Range("E27").Select()
